$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 710, shifting existing rows 710:751 down to 711:752
$ws.Rows.Item(710).Insert()

# Populate the newly inserted row 710 with the new data point.
# Column A holds a date-like string ("2026/01/28") that must stay plain text
# (matching every other row in the sheet) rather than being auto-converted
# into a date serial number, so the cell is temporarily forced to Text format
# while the value is assigned, then returned to the default "Normal" style
# (the other data rows carry no explicit style either).
$ws.Cells.Item(710, 1).NumberFormat = "@"
$ws.Cells.Item(710, 1).Value = "2026/01/28"
$ws.Cells.Item(710, 1).Style = "Normal"
$ws.Cells.Item(710, 2).Value = "水"
$ws.Cells.Item(710, 3).Value = 3
$ws.Cells.Item(710, 4).Value = 31
